$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1319.0741
$ws.Range("I15").Value = 1319.0741
$ws.Range("K15").Value = 3957.2223
$ws.Range("M15").Value = -3788.2223
$ws.Range("H33").Value = 408.79166
$ws.Range("I33").Value = 351.55554
$ws.Range("K33").Value = 351.55554
$ws.Range("M33").Value = -122.55554
$ws.Range("H62").Value = 7983.3335
$ws.Range("I62").Value = 7900
$ws.Range("K62").Value = 7900
$ws.Range("M62").Value = -7276
$ws.Range("H65").Value = 7983.3335
$ws.Range("I65").Value = 7900
$ws.Range("K65").Value = 39500
$ws.Range("M65").Value = -36380
$ws.Range("H86").Value = 6107
$ws.Range("I86").Value = 6096.75
$ws.Range("K86").Value = 6096.75
$ws.Range("M86").Value = -4973.75
$ws.Range("H87").Value = 29875
$ws.Range("H89").Value = 6107
$ws.Range("I89").Value = 6096.75
$ws.Range("K89").Value = 30483.75
$ws.Range("M89").Value = -24867.75
$ws.Range("H90").Value = 29875
$ws.Range("H98").Value = 2590.3845
$ws.Range("I98").Value = 1722.5
$ws.Range("K98").Value = 1722.5
$ws.Range("M98").Value = -224.5
$ws.Range("H101").Value = 3869.125
$ws.Range("I101").Value = 1988.5
$ws.Range("J101").Value = 5749.75
$ws.Range("K101").Value = 5965.5
$ws.Range("L101").Value = 17249.25
$ws.Range("M101").Value = -4343.5
$ws.Range("N101").Value = -20493.25
$ws.Range("H112").Value = 2276.2273
$ws.Range("J112").Value = 2276.2273
$ws.Range("L112").Value = 6828.6819
$ws.Range("N112").Value = -9044.6819
$ws.Range("H115").Value = 390.3846
$ws.Range("I115").Value = 390.3846
$ws.Range("K115").Value = 1171.1538
$ws.Range("M115").Value = 395.8462
$ws.Range("H122").Value = 2590.3845
$ws.Range("I122").Value = 1722.5
$ws.Range("K122").Value = 5167.5
$ws.Range("M122").Value = -2717.5
$ws.Range("H129").Value = 2023.5
$ws.Range("J129").Value = 2874.1667
$ws.Range("L129").Value = 8622.500100000001
$ws.Range("N129").Value = -18622.5001
$ws.Range("H132").Value = 1539.2291
$ws.Range("I132").Value = 1539.2291
$ws.Range("K132").Value = 4617.6873
$ws.Range("M132").Value = -2087.6873
$ws.Range("H137").Value = 11767555
$ws.Range("I137").Value = 47620948
$ws.Range("J137").Value = 3160.7188
$ws.Range("K137").Value = 142862844
$ws.Range("L137").Value = 9482.1564
$ws.Range("M137").Value = -142860294
$ws.Range("N137").Value = -14582.1564
$ws.Range("H138").Value = 4382.3096
$ws.Range("I138").Value = 2430.2144
$ws.Range("J138").Value = 5358.357
$ws.Range("K138").Value = 7290.6432
$ws.Range("L138").Value = 16075.071
$ws.Range("M138").Value = -2150.6432
$ws.Range("N138").Value = -26355.071

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H18").Value = 6250
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()
$ws.Range("H32").Value = 16932.615
$ws.Range("I32").Value = 16932.615
$ws.Range("K32").Value = 16932.615
$ws.Range("M32").Value = -16645.615

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 7563
$ws.Range("I20").Value = 4840
$ws.Range("J20").Value = 13009
$ws.Range("K20").Value = 4840
$ws.Range("L20").Value = 13009
$ws.Range("M20").Value = -4593
$ws.Range("N20").Value = -13503
$ws.Range("H94").Value = 3252.5652
$ws.Range("I94").Value = 3191.7
$ws.Range("J94").Value = 3658.3333
$ws.Range("K94").Value = 3191.7
$ws.Range("L94").Value = 3658.3333
$ws.Range("M94").Value = -2740.7
$ws.Range("N94").Value = -4560.3333
$ws.Range("H134").Value = 2699.6667
$ws.Range("I134").Value = 2013.9259
$ws.Range("K134").Value = 6041.7777
$ws.Range("M134").Value = -3506.7777

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H45").Value = 7500
$ws.Range("I45").Value = 7500
$ws.Range("K45").Value = 7500
$ws.Range("M45").Value = -6907
$ws.Range("H58").Value = 4227.52
$ws.Range("I58").Value = 1310.375
$ws.Range("J58").Value = 5600.294
$ws.Range("K58").Value = 1310.375
$ws.Range("L58").Value = 5600.294
$ws.Range("M58").Value = -1107.375
$ws.Range("N58").Value = -6006.294
$ws.Range("H134").Value = 4948.4165
$ws.Range("I134").Value = 3436.7
$ws.Range("K134").Value = 10310.1
$ws.Range("M134").Value = -7775.099999999999
$ws.Range("H136").Value = 4227.52
$ws.Range("I136").Value = 1310.375
$ws.Range("J136").Value = 5600.294
$ws.Range("K136").Value = 3931.125
$ws.Range("L136").Value = 16800.882
$ws.Range("M136").Value = -1381.125
$ws.Range("N136").Value = -21900.882

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H126").Value = 3190
$ws.Range("I126").Value = 1316.6666
$ws.Range("J126").Value = 6000
$ws.Range("K126").Value = 3949.9998
$ws.Range("L126").Value = 18000
$ws.Range("M126").Value = 990.0001999999999
$ws.Range("N126").Value = -27880

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 15998.777
$ws.Range("I70").Value = 14855.857
$ws.Range("J70").Value = 19999
$ws.Range("K70").Value = 14855.857
$ws.Range("L70").Value = 19999
$ws.Range("M70").Value = -14585.857
$ws.Range("N70").Value = -20539
$ws.Range("H73").Value = 15998.777
$ws.Range("I73").Value = 14855.857
$ws.Range("J73").Value = 19999
$ws.Range("K73").Value = 14855.857
$ws.Range("L73").Value = 19999
$ws.Range("M73").Value = -13919.857
$ws.Range("N73").Value = -21871
$ws.Range("H80").Value = 3996.0417
$ws.Range("I80").Value = 2860.923
$ws.Range("J80").Value = 5337.5454
$ws.Range("K80").Value = 2860.923
$ws.Range("L80").Value = 5337.5454
$ws.Range("M80").Value = -1862.923
$ws.Range("N80").Value = -7333.5454
$ws.Range("H83").Value = 3996.0417
$ws.Range("I83").Value = 2860.923
$ws.Range("J83").Value = 5337.5454
$ws.Range("K83").Value = 14304.615
$ws.Range("L83").Value = 26687.727
$ws.Range("M83").Value = -9312.614999999998
$ws.Range("N83").Value = -36671.727

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6672.7354
$ws.Range("I7").Value = 5537.769
$ws.Range("J7").Value = 10361.375
$ws.Range("K7").Value = 5537.769
$ws.Range("L7").Value = 10361.375
$ws.Range("M7").Value = -5425.769
$ws.Range("N7").Value = -10585.375
$ws.Range("H55").Value = 3291.56
$ws.Range("I55").Value = 1908.3636
$ws.Range("J55").Value = 4378.357
$ws.Range("K55").Value = 1908.3636
$ws.Range("L55").Value = 4378.357
$ws.Range("M55").Value = -1735.3636
$ws.Range("N55").Value = -4724.357
$ws.Range("H61").Value = 6192.8237
$ws.Range("I61").Value = 6201.5
$ws.Range("J61").Value = 6180.4287
$ws.Range("K61").Value = 6201.5
$ws.Range("L61").Value = 6180.4287
$ws.Range("M61").Value = -5999.5
$ws.Range("N61").Value = -6584.4287
$ws.Range("H68").Value = 7897.45
$ws.Range("I68").Value = 4020.8572
$ws.Range("K68").Value = 4020.8572
$ws.Range("M68").Value = -3271.8572
$ws.Range("H71").Value = 7897.45
$ws.Range("I71").Value = 4020.8572
$ws.Range("K71").Value = 20104.286
$ws.Range("M71").Value = -16360.286
$ws.Range("H82").Value = 4147.2285
$ws.Range("I82").Value = 1920.3572
$ws.Range("K82").Value = 1920.3572
$ws.Range("M82").Value = -1559.3572
$ws.Range("H85").Value = 4147.2285
$ws.Range("I85").Value = 1920.3572
$ws.Range("K85").Value = 1920.3572
$ws.Range("M85").Value = -672.3571999999999
$ws.Range("H93").Value = 3843.875
$ws.Range("I93").Value = 3875.1667
$ws.Range("K93").Value = 3875.1667
$ws.Range("M93").Value = -2627.1667
$ws.Range("H113").Value = 6192.8237
$ws.Range("I113").Value = 6201.5
$ws.Range("J113").Value = 6180.4287
$ws.Range("K113").Value = 6201.5
$ws.Range("L113").Value = 6180.4287
$ws.Range("M113").Value = -4031.5
$ws.Range("N113").Value = -10520.4287
$ws.Range("H126").Value = 6672.7354
$ws.Range("I126").Value = 5537.769
$ws.Range("J126").Value = 10361.375
$ws.Range("K126").Value = 16613.307
$ws.Range("L126").Value = 31084.125
$ws.Range("M126").Value = -14143.307
$ws.Range("N126").Value = -36024.125
$ws.Range("H136").Value = 4462.1
$ws.Range("I136").Value = 2847.3845
$ws.Range("J136").Value = 7460.857
$ws.Range("K136").Value = 8542.1535
$ws.Range("L136").Value = 22382.571
$ws.Range("M136").Value = -5992.1535
$ws.Range("N136").Value = -27482.571

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7449.9165
$ws.Range("I62").Value = 6680
$ws.Range("J62").Value = 7999.857
$ws.Range("K62").Value = 6680
$ws.Range("L62").Value = 7999.857
$ws.Range("M62").Value = -6056
$ws.Range("N62").Value = -9247.857
$ws.Range("H65").Value = 7449.9165
$ws.Range("I65").Value = 6680
$ws.Range("J65").Value = 7999.857
$ws.Range("K65").Value = 33400
$ws.Range("L65").Value = 39999.285
$ws.Range("M65").Value = -30280
$ws.Range("N65").Value = -46239.285
$ws.Range("H92").Value = 20000
$ws.Range("J92").Value = 20000
$ws.Range("L92").Value = 20000
$ws.Range("N92").Value = -24992
$ws.Range("H107").Value = 1112.2812
$ws.Range("I107").Value = 916.7083
$ws.Range("K107").Value = 2750.1249
$ws.Range("M107").Value = -830.1248999999998
$ws.Range("H132").Value = 1963.8379
$ws.Range("J132").Value = 19005
$ws.Range("L132").Value = 57015
$ws.Range("N132").Value = -62075
$ws.Range("H136").Value = 2902.4
$ws.Range("J136").Value = 7481.4165
$ws.Range("L136").Value = 22444.2495
$ws.Range("N136").Value = -27544.2495
